$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(462, 44249.53087737269, "Tampere", "36-40 v", "mies", 15, "Työntekijä / palkollinen", 1, "fullstack-ohjelmistokehittä / arkkitehti / pilviveikko", "Etä", 5700, 70000, $true, "", "", 5833.333333333333),
    @(463, 44249.53103743056, "Oulu", "26-30 v", "mies", 7, "Työntekijä / palkollinen", 1, "Backend", "Etä", 3800, 47500, $true, "", "", 3958.333333333333),
    @(464, 44249.53438325231, "PK-Seutu", "26-30 v", "mies", 5, "Työntekijä / palkollinen", 1, "Mobiilikehittäjä", "Toimisto", 4500, 56250, $true, "", "", 4687.5),
    @(465, 44249.5357290625, "Oulu", "26-30 v", "nainen", 5, "Työntekijä / palkollinen", 1, "Web developer", "50/50", 3000, 37500, $false, "", "Kokemusta kokonaisuudessaan 7v, mutta siitä reilut kaksi vuotta lasten kanssa kotona koodaamatta.", 3125),
    @(466, 44249.53759880787, "PK-Seutu", "26-30 v", "mies", 9, "Työntekijä / palkollinen", 1, "Tuotepäällikkö", "Toimisto", 5500, 82500, $true, "", "", 6875),
    @(467, 44249.54394976852, "Tampere", "31-35 v", "mies", 5, "Työntekijä / palkollinen", 1, "Lead front end dev", "Toimisto", 4200, 50000, $true, "", "", 4166.666666666667),
    @(468, 44249.56513866898, "PK-Seutu", "26-30 v", "mies", 0, "Työntekijä / palkollinen", 1, "harjoittelija", "Toimisto", 2200, 27500, $false, "", "", 2291.666666666667),
    @(469, 44249.59106795139, "EU", "31-35 v", "mies", 8, "Työntekijä / palkollinen", 1, "Senior Backend Developer", "Toimisto", 4800, 59000, $false, "", "", 4916.666666666667)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
    $ws.Cells.Item($r, 9).Value = $row[9]
    $ws.Cells.Item($r, 10).Value = $row[10]
    $ws.Cells.Item($r, 11).Value = $row[11]
    $ws.Cells.Item($r, 12).Value = $row[12]
    $ws.Cells.Item($r, 13).Value = $row[13]
    $ws.Cells.Item($r, 14).Value = $row[14]
    $ws.Cells.Item($r, 15).Value = $row[15]
}
